$wb = $excel.ActiveWorkbook

# This script applies a scheduled-runner market-price data refresh to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N) across
# several crafting-leve sheets. Values are plain numeric snapshots (no
# formulas anywhere in the workbook), so each touched cell is written
# directly; cells that no longer have a value are cleared so they drop out
# of the saved XML entirely (matching cells that are genuinely absent).

### Sheet: ALC ###
$ws = $wb.Worksheets.Item("ALC")

# Row 2: H2: 196.33333 -> 137.8; J2: 0 -> 50; L2: 0 -> 50; N2: None -> -276
$ws.Range("H2").Value = 137.8
$ws.Range("J2").Value = 50
$ws.Range("L2").Value = 50
$ws.Range("N2").Value = -276

# Row 40: H40: 2659.9 -> 2627.6667; J40: 2810 -> 2775; L40: 2810 -> 2775; N40: -3160 -> -3125
$ws.Range("H40").Value = 2627.6667
$ws.Range("J40").Value = 2775
$ws.Range("L40").Value = 2775
$ws.Range("N40").Value = -3125

# Row 80: H80: 1040576.94 -> 953874.6; I80: 1757442.2 -> 1523136.5; K80: 5272326.6 -> 4569409.5; M80: -5271328.6 -> -4568411.5
$ws.Range("H80").Value = 953874.6
$ws.Range("I80").Value = 1523136.5
$ws.Range("K80").Value = 4569409.5
$ws.Range("M80").Value = -4568411.5

# Row 83: H83: 1040576.94 -> 953874.6; I83: 1757442.2 -> 1523136.5; K83: 15816979.8 -> 13708228.5; M83: -15811987.8 -> -13703236.5
$ws.Range("H83").Value = 953874.6
$ws.Range("I83").Value = 1523136.5
$ws.Range("K83").Value = 13708228.5
$ws.Range("M83").Value = -13703236.5

# Row 92: H92: 133551.27 -> 117870.234; I92: 1166.091 -> 1027.0769; K92: 1166.091 -> 1027.0769; M92: 81.90900000000011 -> 220.9231
$ws.Range("H92").Value = 117870.234
$ws.Range("I92").Value = 1027.0769
$ws.Range("K92").Value = 1027.0769
$ws.Range("M92").Value = 220.9231

# Row 93: H93: 0 -> 54444; I93: 0 -> 54444; K93: 0 -> 54444; M93: None -> -51948
$ws.Range("H93").Value = 54444
$ws.Range("I93").Value = 54444
$ws.Range("K93").Value = 54444
$ws.Range("M93").Value = -51948

# Row 101: H101: 511.1875 -> 490.17648; I101: 319.54544 -> 305.75; K101: 958.63632 -> 917.25; M101: 663.36368 -> 704.75
$ws.Range("H101").Value = 490.17648
$ws.Range("I101").Value = 305.75
$ws.Range("K101").Value = 917.25
$ws.Range("M101").Value = 704.75

# Row 132: H132: 4418.676 -> 4438.6665; J132: 5606.6924 -> 5765.6665; L132: 16820.0772 -> 17296.9995; N132: -21880.0772 -> -22356.9995
$ws.Range("H132").Value = 4438.6665
$ws.Range("J132").Value = 5765.6665
$ws.Range("L132").Value = 17296.9995
$ws.Range("N132").Value = -22356.9995

# Row 138: H138: 2015.2 -> 1975.9032; I138: 1712 -> 1641.6154; K138: 5136 -> 4924.8462; M138: 4 -> 215.1538
$ws.Range("H138").Value = 1975.9032
$ws.Range("I138").Value = 1641.6154
$ws.Range("K138").Value = 4924.8462
$ws.Range("M138").Value = 215.1538

### Sheet: ARM ###
$ws = $wb.Worksheets.Item("ARM")

# Row 32: H32: 4385.7 -> 4291.5366; I32: 3200.7568 -> 3130.342; K32: 3200.7568 -> 3130.342; M32: -2913.7568 -> -2843.342
$ws.Range("H32").Value = 4291.5366
$ws.Range("I32").Value = 3130.342
$ws.Range("K32").Value = 3130.342
$ws.Range("M32").Value = -2843.342

# Row 46: H46: 6013.8335 -> 5646.857; I46: 4969.5 -> 4461.3335; K46: 4969.5 -> 4461.3335; M46: -4650.5 -> -4142.3335
$ws.Range("H46").Value = 5646.857
$ws.Range("I46").Value = 4461.3335
$ws.Range("K46").Value = 4461.3335
$ws.Range("M46").Value = -4142.3335

# Row 102: H102: 52333.73 -> 52331.27; I102: 45269.258 -> 45266.324; K102: 45269.258 -> 45266.324; M102: -43647.258 -> -43644.324
$ws.Range("H102").Value = 52331.27
$ws.Range("I102").Value = 45266.324
$ws.Range("K102").Value = 45266.324
$ws.Range("M102").Value = -43644.324

### Sheet: BSM ###
$ws = $wb.Worksheets.Item("BSM")

# Row 107: H107: 23811594 -> 23811592; I107: 1284.2858 -> 1275; K107: 1284.2858 -> 1275; M107: 635.7141999999999 -> 645
$ws.Range("H107").Value = 23811592
$ws.Range("I107").Value = 1275
$ws.Range("K107").Value = 1275
$ws.Range("M107").Value = 645

# Row 134: H134: 1591.1455 -> 1616.4073; I134: 1452.8302 -> 1476.4038; K134: 4358.4906 -> 4429.2114; M134: -1823.4906 -> -1894.2114
$ws.Range("H134").Value = 1616.4073
$ws.Range("I134").Value = 1476.4038
$ws.Range("K134").Value = 4429.2114
$ws.Range("M134").Value = -1894.2114

### Sheet: CRP ###
$ws = $wb.Worksheets.Item("CRP")

# Row 86: H86: 125056 -> 125035.29; I86: 211875.5 -> 211839.25; K86: 211875.5 -> 211839.25; M86: -210752.5 -> -210716.25
$ws.Range("H86").Value = 125035.29
$ws.Range("I86").Value = 211839.25
$ws.Range("K86").Value = 211839.25
$ws.Range("M86").Value = -210716.25

# Row 89: H89: 125056 -> 125035.29; I89: 211875.5 -> 211839.25; K89: 1059377.5 -> 1059196.25; M89: -1053761.5 -> -1053580.25
$ws.Range("H89").Value = 125035.29
$ws.Range("I89").Value = 211839.25
$ws.Range("K89").Value = 1059196.25
$ws.Range("M89").Value = -1053580.25

# Row 105: H105: 759.1429000000001 -> 613.625; I105: 759.1429000000001 -> 613.625; K105: 759.1429000000001 -> 613.625; M105: 987.8570999999999 -> 1133.375
$ws.Range("H105").Value = 613.625
$ws.Range("I105").Value = 613.625
$ws.Range("K105").Value = 613.625
$ws.Range("M105").Value = 1133.375

# Row 139: H139: 0 -> 80780; J139: 0 -> 80780; L139: 0 -> 80780; N139: None -> -91060
$ws.Range("H139").Value = 80780
$ws.Range("J139").Value = 80780
$ws.Range("L139").Value = 80780
$ws.Range("N139").Value = -91060

### Sheet: CUL ###
$ws = $wb.Worksheets.Item("CUL")

# Row 35: H35: 192.5 -> 178.5; I35: 150 -> 122; K35: 450 -> 366; M35: -162 -> -78
$ws.Range("H35").Value = 178.5
$ws.Range("I35").Value = 122
$ws.Range("K35").Value = 366
$ws.Range("M35").Value = -78

# Row 81: H81: 6110.4443 -> 5748.8; J81: 9999.75 -> 8498.6; L81: 29999.25 -> 25495.8; N81: -32245.25 -> -27741.8
$ws.Range("H81").Value = 5748.8
$ws.Range("J81").Value = 8498.6
$ws.Range("L81").Value = 25495.8
$ws.Range("N81").Value = -27741.8

# Row 84: H84: 6110.4443 -> 5748.8; J84: 9999.75 -> 8498.6; L84: 89997.75 -> 76487.40000000001; N84: -101229.75 -> -87719.40000000001
$ws.Range("H84").Value = 5748.8
$ws.Range("J84").Value = 8498.6
$ws.Range("L84").Value = 76487.40000000001
$ws.Range("N84").Value = -87719.40000000001

# Row 121: H121: 736138.7 -> 736210; I121: 1133478.9 -> 1275148.1; J121: 20926.4 -> 17625.834; K121: 3400436.7 -> 3825444.3; L121: 62779.2 -> 52877.50199999999; M121: -3399126.7 -> -3824134.3; N121: -65399.2 -> -55497.50199999999
$ws.Range("H121").Value = 736210
$ws.Range("I121").Value = 1275148.1
$ws.Range("J121").Value = 17625.834
$ws.Range("K121").Value = 3825444.3
$ws.Range("L121").Value = 52877.50199999999
$ws.Range("M121").Value = -3824134.3
$ws.Range("N121").Value = -55497.50199999999

### Sheet: GSM ###
$ws = $wb.Worksheets.Item("GSM")

# Row 102: H102: 1789.8334 -> 1789.9166; I102: 886.84 -> 886.96; K102: 886.84 -> 886.96; M102: 735.16 -> 735.04
$ws.Range("H102").Value = 1789.9166
$ws.Range("I102").Value = 886.96
$ws.Range("K102").Value = 886.96
$ws.Range("M102").Value = 735.04

### Sheet: LTW ###
$ws = $wb.Worksheets.Item("LTW")

# Row 3: H3: 20000 -> 0; I3: 20000 -> 0; K3: 20000 -> 0; M3: -19888 -> None
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 7: H7: 20003056 -> 20836418; I7: 29414476 -> 31252734; K7: 29414476 -> 31252734; M7: -29414364 -> -31252622
$ws.Range("H7").Value = 20836418
$ws.Range("I7").Value = 31252734
$ws.Range("K7").Value = 31252734
$ws.Range("M7").Value = -31252622

# Row 9: H9: 1007.5 -> 975; I9: 1007.5 -> 975; K9: 1007.5 -> 975; M9: -783.5 -> -751
$ws.Range("H9").Value = 975
$ws.Range("I9").Value = 975
$ws.Range("K9").Value = 975
$ws.Range("M9").Value = -751

# Row 14: H14: 0 -> 5000; I14: 0 -> 5000; K14: 0 -> 5000; M14: None -> -4828
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4828

# Row 15: H15: 20000 -> 0; I15: 20000 -> 0; K15: 20000 -> 0; M15: -19830 -> None
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# Row 40: H40: 3119.125 -> 6383.3335; I40: 2652.8096 -> 0; K40: 2652.8096 -> 0; M40: -2516.8096 -> None
$ws.Range("H40").Value = 6383.3335
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 61: H61: 2153.1 -> 2184.5454; J61: 3750 -> 3333; L61: 3750 -> 3333; N61: -4154 -> -3737
$ws.Range("H61").Value = 2184.5454
$ws.Range("J61").Value = 3333
$ws.Range("L61").Value = 3333
$ws.Range("N61").Value = -3737

# Row 93: H93: 18522818 -> 15877061; I93: 25645056 -> 19611666; J93: 4995.6 -> 4988.25; K93: 25645056 -> 19611666; L93: 4995.6 -> 4988.25; M93: -25643808 -> -19610418; N93: -7491.6 -> -7484.25
$ws.Range("H93").Value = 15877061
$ws.Range("I93").Value = 19611666
$ws.Range("J93").Value = 4988.25
$ws.Range("K93").Value = 19611666
$ws.Range("L93").Value = 4988.25
$ws.Range("M93").Value = -19610418
$ws.Range("N93").Value = -7484.25

# Row 100: H100: 2638.5217 -> 2735.7273; I100: 2594.05 -> 2704.2632; K100: 2594.05 -> 2704.2632; M100: -2053.05 -> -2163.2632
$ws.Range("H100").Value = 2735.7273
$ws.Range("I100").Value = 2704.2632
$ws.Range("K100").Value = 2704.2632
$ws.Range("M100").Value = -2163.2632

# Row 113: H113: 2153.1 -> 2184.5454; J113: 3750 -> 3333; L113: 3750 -> 3333; N113: -8090 -> -7673
$ws.Range("H113").Value = 2184.5454
$ws.Range("J113").Value = 3333
$ws.Range("L113").Value = 3333
$ws.Range("N113").Value = -7673

# Row 122: H122: 4366.3335 -> 4799.5; I122: 3549.5 -> 3599; K122: 10648.5 -> 10797; M122: -8198.5 -> -8347
$ws.Range("H122").Value = 4799.5
$ws.Range("I122").Value = 3599
$ws.Range("K122").Value = 10797
$ws.Range("M122").Value = -8347

# Row 126: H126: 20003056 -> 20836418; I126: 29414476 -> 31252734; K126: 88243428 -> 93758202; M126: -88240958 -> -93755732
$ws.Range("H126").Value = 20836418
$ws.Range("I126").Value = 31252734
$ws.Range("K126").Value = 93758202
$ws.Range("M126").Value = -93755732

# Row 132: H132: 41673590 -> 90920370; I132: 71431780 -> 333336260; J132: 12127.4 -> 14409.25; K132: 214295340 -> 1000008780; L132: 36382.2 -> 43227.75; M132: -214292810 -> -1000006250; N132: -41442.2 -> -48287.75
$ws.Range("H132").Value = 90920370
$ws.Range("I132").Value = 333336260
$ws.Range("J132").Value = 14409.25
$ws.Range("K132").Value = 1000008780
$ws.Range("L132").Value = 43227.75
$ws.Range("M132").Value = -1000006250
$ws.Range("N132").Value = -48287.75

### Sheet: WVR ###
$ws = $wb.Worksheets.Item("WVR")

# Row 9: H9: 0 -> 5000; I9: 0 -> 5000; K9: 0 -> 5000; M9: None -> -4860
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 5000
$ws.Range("K9").Value = 5000
$ws.Range("M9").Value = -4860

# Row 31: H31: 35000 -> 0; J31: 35000 -> 0; L31: 35000 -> 0; N31: -35696 -> None
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

# Row 122: H122: 1965.5853 -> 2040.3846; I122: 1782 -> 1799.2858; J122: 2361 -> 2654.0908; K122: 5346 -> 5397.857400000001; L122: 7083 -> 7962.2724; M122: -2896 -> -2947.857400000001; N122: -11983 -> -12862.2724
$ws.Range("H122").Value = 2040.3846
$ws.Range("I122").Value = 1799.2858
$ws.Range("J122").Value = 2654.0908
$ws.Range("K122").Value = 5397.857400000001
$ws.Range("L122").Value = 7962.2724
$ws.Range("M122").Value = -2947.857400000001
$ws.Range("N122").Value = -12862.2724

# Row 126: H126: 1853.5 -> 1729.909; J126: 2486.5 -> 2088; L126: 7459.5 -> 6264; N126: -12399.5 -> -11204
$ws.Range("H126").Value = 1729.909
$ws.Range("J126").Value = 2088
$ws.Range("L126").Value = 6264
$ws.Range("N126").Value = -11204

# Row 132: H132: 8231.25 -> 7674.923; J132: 7000 -> 5799.8; L132: 21000 -> 17399.4; N132: -26060 -> -22459.4
$ws.Range("H132").Value = 7674.923
$ws.Range("J132").Value = 5799.8
$ws.Range("L132").Value = 17399.4
$ws.Range("N132").Value = -22459.4
